$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DATA")
$ws.Range("B3:L3").ClearContents()
